# Skin Cancer Final.pptx - "Add files via upload" edit
#
# Semantic changes reproduced here:
#   1. Slide 14 ("Final / Summary" slide), the big bullet text box:
#        - first sentence rewritten
#        - box height grown (autofit) from 8262000 EMU to 9109633 EMU
#          (== 650.5511811.. pt -> 717.293937.. pt)

$p = $ppt.ActivePresentation

$slide14 = $p.Slides.Item(14)
$summaryShape = $slide14.Shapes.Item(3)

# --- text edit -------------------------------------------------------
$oldSentence = "We found that age is not a true determinant in skin cancer. "
$newSentence = "Sex has a statistically significant effect on the mean age at which they were diagnosed with skin cancer."

$tr = $summaryShape.TextFrame.TextRange
$firstRun = $tr.Characters(1, $oldSentence.Length)
$firstRun.Text = $newSentence

# --- resize (the textbox auto-grew to fit the new wording) -----------
$summaryShape.Height = 9109633 / 12700
